# Applies the changes described in the commit:
#  - CredCard sheet: clear out a block of "reconciliation" helper cells
#    (columns F,G,H,J for several rows; F1:G1 on row 1) while keeping
#    their existing cell styles, and update the selection accordingly.
#  - TestRecord sheet: bump a transaction's date and amount.
#  - Budget Out sheet: bump an amount.
#  - Expected Out sheet: bump two summary amounts (the SUM formula at
#    the top of the sheet recalculates automatically).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# CredCard sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CredCard")
$ws.Activate()

$ws.Range("F1:G1").ClearContents()

$ws.Range("F3:H3").ClearContents()
$ws.Range("J3").ClearContents()

$ws.Range("F4:H4").ClearContents()
$ws.Range("J4").ClearContents()

$ws.Range("F6:H6").ClearContents()
$ws.Range("J6").ClearContents()

$ws.Range("F7:H7").ClearContents()
$ws.Range("J7").ClearContents()

$ws.Range("F8:H8").ClearContents()
$ws.Range("J8").ClearContents()

$ws.Range("F10:H10").ClearContents()
$ws.Range("J10").ClearContents()

$ws.Range("F11:H11").ClearContents()
$ws.Range("J11").ClearContents()

# Update the selection to match the new state (F1:J11, active cell F1)
$ws.Range("F1:J11").Select()

# ---------------------------------------------------------------
# Budget Out sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Budget Out")
$ws.Range("C9").Value = 90.22
$ws.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# ---------------------------------------------------------------
# TestRecord sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TestRecord")
$ws.Range("A10").Value = 43263
$ws.Range("B10").Value = 121.14
$ws.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# ---------------------------------------------------------------
# Expected Out sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Expected Out")
$ws.Range("B9").Value = 1348.16
$ws.Range("B11").Value = 428.02

# Restore CredCard as the active sheet (it is the originally active /
# tab-selected sheet in the workbook).
$wb.Worksheets.Item("CredCard").Activate()
